# Update NATMI TPM-derived statistics for the Tnf-Tnfrsf21 LR-pair sheet
# (new values recomputed with the updated TPM matrix). Only columns
# G,H,I,J (ligand expression / specificity), M,N,O,P (receptor expression /
# specificity) and Q,R,S,T (edge weight / specificity) change; columns
# A-F, K, L are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 26
$rowCount = $lastDataRow - $firstDataRow + 1

# Ligand-side values (cols G:J) repeat per 5-row "sending cluster" block
# (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac).
$ligandBlocks = @(
    @([double]"4.300877666666666", [double]"12.902633", [double]"0.04061703229494078", [double]"0.04061703229494078"),
    @([double]"0.3544293333333333", [double]"1.063288", [double]"0.003347193013613811", [double]"0.003347193013613811"),
    @([double]"61.65203333333333", [double]"184.9561", [double]"0.582235260574047", [double]"0.5822352605740471"),
    @([double]"0.10468", [double]"0.31404", [double]"0.0009885868118471018", [double]"0.0009885868118471018"),
    @([double]"39.47650533333334", [double]"118.429516", [double]"0.3728119273055513", [double]"0.3728119273055513")
)

# Receptor-side values (cols M:P) repeat per "target cluster" position
# (ECs, FAPs, Inflammatory-Mac, MuSCs, Resolving-Mac); the same 5 values
# cycle within every sending-cluster block.
$receptorCycle = @(
    @([double]"5.455768666666667", [double]"16.367306", [double]"0.0824390136851795", [double]"0.0824390136851795"),
    @([double]"2.300815", [double]"6.902445", [double]"0.03476630532942922", [double]"0.03476630532942922"),
    @([double]"29.166511", [double]"87.499533", [double]"0.4407185396566677", [double]"0.4407185396566677"),
    @([double]"3.497096", [double]"10.491288", [double]"0.05284262633124592", [double]"0.05284262633124593"),
    @([double]"25.75926033333333", [double]"77.277781", [double]"0.3892335149974776", [double]"0.3892335149974776")
)

# Edge weight / specificity values (cols Q:T) are unique for every data row.
$edgeValues = @(
    @([double]"23.46459361296644", [double]"211.181342516698", [double]"0.003348428081214", [double]"0.003348428081214001"),
    @([double]"9.895523848631665", [double]"89.05971463768499", [double]"0.001412104146341198", [double]"0.001412104146341199"),
    @([double]"125.4415957744876", [double]"1128.974361970389", [double]"0.01790067915821401", [double]"0.01790067915821401"),
    @([double]"15.04058208458933", [double]"135.365238761304", [double]"0.002146310660245703", [double]"0.002146310660245704"),
    @([double]"110.7874274774859", [double]"997.0868472973729", [double]"0.01580951024892586", [double]"0.01580951024892586"),
    @([double]"1.933684451347556", [double]"17.403160062128", [double]"0.0002759392906562461", [double]"0.0002759392906562461"),
    @([double]"0.8154763265733334", [double]"7.33928693916", [double]"0.0001163695343078301", [double]"0.0001163695343078301"),
    @([double]"10.33746704938933", [double]"93.037203444504", [double]"0.00147517001690888", [double]"0.00147517001690888"),
    @([double]"1.239473403882666", [double]"11.155260634944", [double]"0.0001768744696769515", [double]"0.0001768744696769516"),
    @([double]"9.12983746710311", [double]"82.16853720392801", [double]"0.001302839702063903", [double]"0.001302839702063903"),
    @([double]"336.3592316962889", [double]"3027.2330852666", [double]"0.04799890061445791", [double]"0.04799890061445792"),
    @([double]"141.8499230738333", [double]"1276.6493076645", [double]"0.0202421688426771", [double]"0.02024216884267711"),
    @([double]"1798.174708389033", [double]"16183.5723755013", [double]"0.2566018737768134", [double]"0.2566018737768134"),
    @([double]"215.6030791618666", [double]"1940.4277124568", [double]"0.03076684031138997", [double]"0.03076684031138998"),
    @([double]"1588.110776712678", [double]"14292.9969904141", [double]"0.2266254770287086", [double]"0.2266254770287086"),
    @([double]"0.5711098640266666", [double]"5.139988776239999", [double]"8.149812171085119E-05", [double]"8.149812171085119E-05"),
    @([double]"0.2408493142", [double]"2.1676438278", [double]"3.436951094532334E-05", [double]"3.436951094532334E-05"),
    @([double]"3.05315037148", [double]"27.47835334332", [double]"0.0004356885360410957", [double]"0.0004356885360410957"),
    @([double]"0.3660760092799999", [double]"3.29468408352", [double]"5.223952349443412E-05", [double]"5.223952349443413E-05"),
    @([double]"2.696479371693333", [double]"24.26831434524", [double]"0.0003847911196553975", [double]"0.0003847911196553975"),
    @([double]"215.3746808670996", [double]"1938.372127803896", [double]"0.03073424757714049", [double]"0.03073424757714049"),
    @([double]"90.82813561851334", [double]"817.4532205666201", [double]"0.01296129329515777", [double]"0.01296129329515777"),
    @([double]"1151.391927046225", [double]"10362.52734341603", [double]"0.1643051281686903", [double]"0.1643051281686903"),
    @([double]"138.0531288951786", [double]"1242.478160056608", [double]"0.01970036136643886", [double]"0.01970036136643887"),
    @([double]"1016.885577931555", [double]"9151.970201383998", [double]"0.1451108968981238", [double]"0.1451108968981238")
)

$ligandArr = New-Object "object[,]" $rowCount,4
$receptorArr = New-Object "object[,]" $rowCount,4
$edgeArr = New-Object "object[,]" $rowCount,4

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $i = $row - $firstDataRow
    $blockIndex = [math]::Floor($i / 5)
    $cycleIndex = $i % 5

    $ligand = $ligandBlocks[$blockIndex]
    $receptor = $receptorCycle[$cycleIndex]
    $edge = $edgeValues[$i]

    for ($c = 0; $c -lt 4; $c++) {
        $ligandArr[$i, $c] = $ligand[$c]
        $receptorArr[$i, $c] = $receptor[$c]
        $edgeArr[$i, $c] = $edge[$c]
    }
}

$ws.Range("G$($firstDataRow):J$($lastDataRow)").Value = $ligandArr
$ws.Range("M$($firstDataRow):P$($lastDataRow)").Value = $receptorArr
$ws.Range("Q$($firstDataRow):T$($lastDataRow)").Value = $edgeArr

